$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the "Acelga" weekly block (old rows 468-557
# shift down to 471-560), making room for a new week of data (Fecha 44504).
$ws.Rows("468:470").Insert()

# --- New row 468: Calidad "Extra" ---
$ws.Cells.Item(468,1).Value  = 6
$ws.Cells.Item(468,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(468,3).Value  = "Metropolitana"
$ws.Cells.Item(468,4).Value  = 44504
$ws.Cells.Item(468,5).Value  = 13
$ws.Cells.Item(468,6).Value  = 100112009
$ws.Cells.Item(468,7).Value  = "Acelga"
$ws.Cells.Item(468,8).Value  = "Sin especificar"
$ws.Cells.Item(468,9).Value  = "Extra"
$ws.Cells.Item(468,10).Value = 110
$ws.Cells.Item(468,11).Value = 12000
$ws.Cells.Item(468,12).Value = 12000
$ws.Cells.Item(468,13).Value = 12000
$ws.Cells.Item(468,14).Value = "$/docena de atados"
$ws.Cells.Item(468,15).Value = "Región Metropolitana"
$ws.Cells.Item(468,16).Value = 4000
$ws.Cells.Item(468,17).Value = 3
$ws.Cells.Item(468,18).Value = "Hortaliza"

# --- New row 469: Calidad "Primera" ---
$ws.Cells.Item(469,1).Value  = 6
$ws.Cells.Item(469,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(469,3).Value  = "Metropolitana"
$ws.Cells.Item(469,4).Value  = 44504
$ws.Cells.Item(469,5).Value  = 13
$ws.Cells.Item(469,6).Value  = 100112009
$ws.Cells.Item(469,7).Value  = "Acelga"
$ws.Cells.Item(469,8).Value  = "Sin especificar"
$ws.Cells.Item(469,9).Value  = "Primera"
$ws.Cells.Item(469,10).Value = 160
$ws.Cells.Item(469,11).Value = 10000
$ws.Cells.Item(469,12).Value = 10000
$ws.Cells.Item(469,13).Value = 10000
$ws.Cells.Item(469,14).Value = "$/docena de atados"
$ws.Cells.Item(469,15).Value = "Región Metropolitana"
$ws.Cells.Item(469,16).Value = 3333
$ws.Cells.Item(469,17).Value = 3
$ws.Cells.Item(469,18).Value = "Hortaliza"

# --- New row 470: Calidad "Segunda" ---
$ws.Cells.Item(470,1).Value  = 6
$ws.Cells.Item(470,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(470,3).Value  = "Metropolitana"
$ws.Cells.Item(470,4).Value  = 44504
$ws.Cells.Item(470,5).Value  = 13
$ws.Cells.Item(470,6).Value  = 100112009
$ws.Cells.Item(470,7).Value  = "Acelga"
$ws.Cells.Item(470,8).Value  = "Sin especificar"
$ws.Cells.Item(470,9).Value  = "Segunda"
$ws.Cells.Item(470,10).Value = 130
$ws.Cells.Item(470,11).Value = 8000
$ws.Cells.Item(470,12).Value = 8000
$ws.Cells.Item(470,13).Value = 8000
$ws.Cells.Item(470,14).Value = "$/docena de atados"
$ws.Cells.Item(470,15).Value = "Región Metropolitana"
$ws.Cells.Item(470,16).Value = 2667
$ws.Cells.Item(470,17).Value = 3
$ws.Cells.Item(470,18).Value = "Hortaliza"
